# Completed the code for Facility and Practitioner Get API
#
# The sheet "ChaseData" holds a small RGID / Provider lookup table:
#   A1: RGID       B1: Prvoider
#   A2: RG-34565541998716   B2: P-535468020225
#
# Update the provider id in B2 and move the active selection to B3
# (matching the saved cursor position recorded in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChaseData")

# Update the provider value.
$ws.Range("B2").Value = "F-147747020225"

# Leave the selection where the author last left it when saving.
$ws.Range("B3").Select()
